$p = $ppt.ActivePresentation
Write-Host "Slides: $($p.Slides.Count)"
